# Insert a new weekly price record as row 65 ("Cebollín" - Vega Monumental
# Concepción), shifting all subsequent rows (old 65-117) down by one
# (to 66-118), matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 65..117 down to 66..118, carrying formatting along (this is
# also what updates the sheet's dimension from A1:R117 to A1:R118).
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new record's data.
$ws.Range("A65").Value = 11
$ws.Range("B65").Value = "Vega Monumental Concepción"
$ws.Range("C65").Value = "Bíobío"
$ws.Range("D65").Value = 45090
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = 100112037
$ws.Range("G65").Value = "Cebollín"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 3200
$ws.Range("L65").Value = 3500
$ws.Range("M65").Value = 3350
$ws.Range("N65").Value = "`$/paquete 36 unidades"
$ws.Range("O65").Value = "Región Metropolitana"
$ws.Range("P65").Value = 93
$ws.Range("Q65").Value = 36
$ws.Range("R65").Value = "Hortaliza"
